$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 1
    7 = 0
    8 = 1
    9 = 0
    10 = 0
    11 = 0
    12 = 2
    13 = 3
    14 = 1
    15 = 0
    16 = 4
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 5
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 1
    39 = 0
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 1
    49 = 0
    50 = 0
    51 = 0
    52 = 1
    53 = 1
    54 = 1
    55 = 0
    56 = 0
    57 = 1
    58 = 2
    59 = 1
    60 = 2
    61 = 1
    62 = 1
    63 = 0
    64 = 1
    65 = 0
    66 = 1
    67 = 1
    68 = 1
    69 = 1
    70 = 0
    71 = 2
    72 = 0
    73 = 0
    74 = 2
    75 = 1
    76 = 2
    77 = 0
    78 = 0
    79 = 1
    80 = 3
    81 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
